$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$insertRow = 38

# Remember existing hyperlinks (row, col, address) before the insert, since
# this engine does not auto-shift Hyperlink.Range when rows are inserted.
$hlInfo = @()
foreach ($hl in $ws.Hyperlinks) {
    $hlInfo += ,@($hl.Range.Row, $hl.Range.Column, $hl.Address)
}

# Insert a new row before row 38, shifting existing rows 38-50 down to 39-51
$ws.Rows.Item($insertRow).Insert()

# Drop the old hyperlinks and re-add them, shifting any at/after the insertion
# point down by one row so they keep pointing at the same (moved) cell.
$ws.Hyperlinks.Delete()
foreach ($info in $hlInfo) {
    $row = $info[0]
    $col = $info[1]
    $addr = $info[2]
    if ($row -ge $insertRow) {
        $row = $row + 1
    }
    $cell = $ws.Cells.Item($row, $col)
    $ws.Hyperlinks.Add($cell, $addr)
}

# Populate the newly inserted row
$ws.Range("A38").Value = "Request to change child support"
$ws.Range("B38").Value = "https://www.illinoislegalaid.org/legal-information/child-support-modification"

# Add the hyperlink for the new row
$ws.Hyperlinks.Add($ws.Range("B38"), "https://www.illinoislegalaid.org/legal-information/child-support-modification")

# Re-adding hyperlinks can introduce a duplicate "Hyperlink" style variant;
# normalize every hyperlinked B cell back onto the sheet's original style.
foreach ($hl in $ws.Hyperlinks) {
    $hl.Range.Style = $ws.Range("B3").Style
}

# Match the final cursor position recorded in the saved workbook.
[void]$ws.Range("B54").Select()
